$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Stable accommodation.jamais.sdf"
$ws.Range("C1").Value = "Unstable accommodation and/or homeless.jamais.sdf"
$ws.Range("D1").Value = "In detention.jamais.sdf"
$ws.Range("E1").Value = "Other.jamais.sdf"
$ws.Range("F1").Value = "Not known / missing.jamais.sdf"
$ws.Range("G1").Value = "Total.jamais.sdf"
